$d = $word.ActiveDocument

$replacements = @(
    @{old="931÷5=186, 1"; new="470÷6=78, 2"},
    @{old="384÷9=42, 6"; new="569÷9=63, 2"},
    @{old="119÷5=23, 4"; new="142÷5=28, 2"},
    @{old="880÷9=97, 7"; new="424÷7=60, 4"},
    @{old="423÷6=70, 3"; new="656÷2=328, 0"},
    @{old="553÷8=69, 1"; new="151÷8=18, 7"},
    @{old="284÷3=94, 2"; new="555÷5=111, 0"},
    @{old="939÷9=104, 3"; new="880÷8=110, 0"},
    @{old="156÷6=26, 0"; new="371÷7=53, 0"},
    @{old="352÷7=50, 2"; new="358÷6=59, 4"},
    @{old="833÷3=277, 2"; new="139÷3=46, 1"},
    @{old="970÷9=107, 7"; new="235÷2=117, 1"},
    @{old="636÷7=90, 6"; new="936÷8=117, 0"},
    @{old="449÷9=49, 8"; new="729÷2=364, 1"},
    @{old="746÷7=106, 4"; new="549÷6=91, 3"},
    @{old="187÷5=37, 2"; new="223÷6=37, 1"},
    @{old="644÷6=107, 2"; new="441÷8=55, 1"},
    @{old="682÷9=75, 7"; new="817÷5=163, 2"},
    @{old="337÷7=48, 1"; new="990÷2=495, 0"},
    @{old="905÷8=113, 1"; new="419÷5=83, 4"},
    @{old="292÷8=36, 4"; new="513÷3=171, 0"},
    @{old="831÷6=138, 3"; new="636÷3=212, 0"},
    @{old="989÷9=109, 8"; new="637÷6=106, 1"},
    @{old="694÷6=115, 4"; new="375÷5=75, 0"},
    @{old="716÷9=79, 5"; new="142÷7=20, 2"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
